# Append the new EUR->ARS quote row (2025-09-17 21:20:31) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a plain "YYYY-MM-DD" text label. Prefix with an apostrophe
# so Excel stores it as literal text instead of auto-converting it to a date.
$ws.Range("A24").Value = "'2025-09-17"
$ws.Range("B24").Value = "21:20:31"
$ws.Range("C24").Value = "1.00 EUR = 1,749.6249"
